$d = $word.ActiveDocument
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Team Name: ______________________________</w:t>
      </w:r>
      <w:r>
        <w:t>____</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Your</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> name: ___________</w:t>
      </w:r>
      <w:r>
        <w:t>__________________</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Poster:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>0</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> – </w:t>
      </w:r>
      <w:r>
        <w:t>abysmal</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> – bad</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> - inadequate</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> - good</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> – great</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>5</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> - exceptional</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Presentation</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">0 – </w:t>
      </w:r>
      <w:r>
        <w:t>abysmal</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>1 – bad</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>2 - inadequate</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>3 - good</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>4 – great</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>5 - exceptional</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>Comments:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Team Name: __________________________________</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Your</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> name: _____________________________</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Quality of Poster and Analysis:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>0</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>5</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>Comments:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Team Name: __________________________________</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Your</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> name: _____________________________</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Quality of Poster and Analysis:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>0</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>5</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>Comments:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Team Name: __________________________________</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Your</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> name: _____________________________</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Quality of Poster and Analysis:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>0</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>5</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>Comments:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>Team Name: ______</w:t>
      </w:r>
      <w:r>
        <w:t>____________________________</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Your</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> name</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:t>___________</w:t>
      </w:r>
      <w:r>
        <w:softHyphen/>
      </w:r>
      <w:r>
        <w:softHyphen/>
      </w:r>
      <w:r>
        <w:softHyphen/>
      </w:r>
      <w:r>
        <w:softHyphen/>
      </w:r>
      <w:r>
        <w:softHyphen/>
        <w:t>__________________</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Quality of Poster and Analysis:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
      </w:pPr>
      <w:r>
        <w:t>0</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>4</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>5</w:t>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0"/>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Comments:</w:t>
      </w:r>
    </w:p>
'@
$d.Content.InsertXML($xml)
Write-Host "Paragraph count:" $d.Paragraphs.Count
